# Apply the "剩余" (remaining) column decrement and the row-95 refill update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row is 99 (header is row 1, data rows 2..99).
$lastRow = 99

for ($r = 2; $r -le $lastRow; $r++) {
    # Row 36 has a malformed start date (F36 = 202510929, not a valid
    # 8-digit yyyymmdd value) so it is left untouched by the daily update.
    if ($r -eq 36) { continue }

    $eCell = $ws.Cells.Item($r, 5)   # column E = "剩余" (remaining)
    $val = $eCell.Value2

    if ($null -eq $val) { continue }

    if ($r -eq 95) {
        # Row 95's remaining count had reached 1 (about to run out), so it
        # gets refilled back to the full total (column D) with a new start
        # date (F) instead of simply decrementing.
        $ws.Cells.Item($r, 5).Value = 10
        $ws.Cells.Item($r, 6).Value = 20251022
    } else {
        $ws.Cells.Item($r, 5).Value = $val - 1
    }
}
